$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.689.92"
$ws.Range("E2").Value = "  -0.77%  "

$ws.Range("D3").Value = "2.308.49"
$ws.Range("E3").Value = "  +2.78%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "270.58"
$ws.Range("E5").Value = "  +0.63%  "

$ws.Range("D6").Value = "93.27"
$ws.Range("E6").Value = "  +6.41%  "

$ws.Range("D7").Value = "0.627"
$ws.Range("E7").Value = "  +0.25%  "

$ws.Range("D9").Value = "0.621"
$ws.Range("E9").Value = "  +0.89%  "

$ws.Range("D10").Value = "44.77"
$ws.Range("E10").Value = "  -2.90%  "

$ws.Range("E11").Value = "  +0.80%  "

$ws.Range("D12").Value = "8.12"
$ws.Range("E12").Value = "  +7.39%  "

$ws.Range("E13").Value = "  +0.19%  "

$ws.Range("D14").Value = "2.650.07"
$ws.Range("E14").Value = "  +2.68%  "

$ws.Range("D15").Value = "15.26"
$ws.Range("E15").Value = "  +1.94%  "

$ws.Range("D16").Value = "0.851"
$ws.Range("E16").Value = "  +6.58%  "

$ws.Range("D17").Value = "2.321.52"
$ws.Range("E17").Value = "  +3.52%  "

$ws.Range("D18").Value = "43.719.46"
$ws.Range("E18").Value = "  -0.66%  "

$ws.Range("E19").Value = "  +1.67%  "

$ws.Range("E20").Value = "  +3.60%  "

$ws.Range("D21").Value = "71.56"
$ws.Range("E21").Value = "  +1.83%  "

$ws.Range("D22").Value = "239.18"
$ws.Range("E22").Value = "  +2.44%  "

$ws.Range("E23").Value = "  -5.06%  "

$ws.Range("D24").Value = "9.68"
$ws.Range("E24").Value = "  +8.43%  "

$ws.Range("E25").Value = "  +0.03%  "

$ws.Range("D26").Value = "11.37"
$ws.Range("E26").Value = "  +3.15%  "

$ws.Range("E27").Value = "  -3.30%  "

$ws.Range("E28").Value = "  +5.42%  "

$ws.Range("E29").Value = "  -4.88%  "

$ws.Range("D30").Value = "39.03"
$ws.Range("E30").Value = "  -4.24%  "

$ws.Range("D31").Value = "22.57"
$ws.Range("E31").Value = "  +8.61%  "

$ws.Range("D32").Value = "171.79"
$ws.Range("E32").Value = "  -2.12%  "

$ws.Range("D33").Value = "0.0901"
$ws.Range("E33").Value = "  -1.17%  "

$ws.Range("D34").Value = "5.58"
$ws.Range("E34").Value = "  +2.75%  "

$ws.Range("E35").Value = "  +1.76%  "

$ws.Range("E36").Value = "  -0.77%  "

$ws.Range("D37").Value = "4.48"
$ws.Range("E37").Value = "  +2.01%  "

$ws.Range("D38").Value = "0.0356"
$ws.Range("E38").Value = "  -0.45%  "

$ws.Range("D39").Value = "3.43"
$ws.Range("E39").Value = "  +2.89%  "

$ws.Range("D40").Value = "0.234"
$ws.Range("E40").Value = "  +14.44%  "

$ws.Range("D41").Value = "2.30"
$ws.Range("E41").Value = "  +6.60%  "

$ws.Range("D42").Value = "12.15"
$ws.Range("E42").Value = "  -4.23%  "

$ws.Range("D43").Value = "1.33"
$ws.Range("E43").Value = "  +16.93%  "

$ws.Range("E44").Value = "  +0.77%  "

$ws.Range("D45").Value = "61.80"
$ws.Range("E45").Value = "  -5.83%  "

$ws.Range("D46").Value = "8.92"
$ws.Range("E46").Value = "  +6.47%  "

$ws.Range("E47").Value = "  +2.77%  "

$ws.Range("D48").Value = "100.24"
$ws.Range("E48").Value = "  -0.12%  "

$ws.Range("E49").Value = "  -1.54%  "

$ws.Range("D50").Value = "2.528.84"
$ws.Range("E50").Value = "  +2.58%  "

$ws.Range("D51").Value = "0.427"
$ws.Range("E51").Value = "  -3.60%  "
